# Adds the mouse processing pipeline rows (Alencar et al.) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Values are entered in this particular order so the shared-string table
# ends up built in the same sequence as the authored workbook.
$ws.Range("A7").Value = "Alencar et al."
$ws.Range("E7").Value = "This is from Gabe Alencar: mouse genome was converted to human, genes were removed (non 1:1 homologues), and then all the libraries were integrated, etc"
$ws.Range("D8").Value = "Aencar_2020_mouse.rds"
$ws.Range("B7").Value = "Alencar_2020_humanized_homolog"
$ws.Range("B8").Value = "Alencar_2020_mouse"
$ws.Range("D7").Value = "Aencar_2020_humanized.rds"
$ws.Range("E8").Value = "This is from Katie Owsiany: unconverted original mouse data"

$ws.Range("C7").Value = "UNPROCESSED.rds"
$ws.Range("C8").Value = "UNPROCESSED.rds"
$ws.Range("A8").Value = "Alencar et al."

# Match formatting style (font size 11) used by existing Authors/DataID columns
$ws.Range("A7:B8").Font.Size = 11

# Widen column B to fit the new longer DataID values (saved width="22")
$ws.Columns.Item(2).ColumnWidth = 21.17

# Update selection to mirror the saved workbook state
$ws.Range("E9").Select()
